# Rename sheets: Sheet1 -> "Input Data", Sheet2 -> "search"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Input Data"
$ws2.Name = "search"

# --- "Input Data" sheet: an email/password row ---
# Write in the same order the source workbook was authored (password1 first,
# so it lands at shared-string index 0, matching Email/password /address/postcode/distance).
$ws1.Range("B2").Value = "password1"
$ws1.Range("A1").Value = "Email"
$ws1.Range("B1").Value = "password "
$ws1.Range("A2").Value = "PrickHead15@yahoo.co.uk"

# Turn the email address into a mailto hyperlink (creates the Hyperlink style too)
$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:PrickHead15@yahoo.co.uk")

# Highlight the header row in yellow
$ws1.Range("A1:B1").Interior.Color = 65535

# Column widths
$ws1.Columns.Item(1).ColumnWidth = 24.140625
$ws1.Columns.Item(2).ColumnWidth = 10.85546875

# --- "search" sheet: postcode/distance headers ---
$ws2.Range("A1").Value = "postcode"
$ws2.Range("B1").Value = "distance"

# Selections: A2 active on "Input Data", B1 active + tab selected on "search"
$null = $ws1.Range("A2").Select()
$null = $ws2.Activate()
$null = $ws2.Range("B1").Select()
